$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5288.8237
$ws.Range("I76").Value = 5767.1665
$ws.Range("J76").Value = 4140.8
$ws.Range("K76").Value = 5767.1665
$ws.Range("L76").Value = 4140.8
$ws.Range("M76").Value = -5452.1665
$ws.Range("N76").Value = -4770.8

$ws.Range("H79").Value = 5288.8237
$ws.Range("I79").Value = 5767.1665
$ws.Range("J79").Value = 4140.8
$ws.Range("K79").Value = 5767.1665
$ws.Range("L79").Value = 4140.8
$ws.Range("M79").Value = -4675.1665
$ws.Range("N79").Value = -6324.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1726.3334
$ws.Range("I2").Value = 1631.625
$ws.Range("J2").Value = 1915.75
$ws.Range("K2").Value = 1631.625
$ws.Range("L2").Value = 1915.75
$ws.Range("M2").Value = -1518.625
$ws.Range("N2").Value = -2141.75

$ws.Range("H116").Value = 1726.3334
$ws.Range("I116").Value = 1631.625
$ws.Range("J116").Value = 1915.75
$ws.Range("K116").Value = 1631.625
$ws.Range("L116").Value = 1915.75
$ws.Range("M116").Value = 662.375
$ws.Range("N116").Value = -6503.75

$ws.Range("H123").Value = 27802.75
$ws.Range("J123").Value = 27802.75
$ws.Range("L123").Value = 27802.75
$ws.Range("N123").Value = -37602.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1726.3334
$ws.Range("I3").Value = 1631.625
$ws.Range("J3").Value = 1915.75
$ws.Range("K3").Value = 1631.625
$ws.Range("L3").Value = 1915.75
$ws.Range("M3").Value = -1517.625
$ws.Range("N3").Value = -2143.75

$ws.Range("H86").Value = 47605.363
$ws.Range("I86").Value = 2332.3125
$ws.Range("J86").Value = 168333.5
$ws.Range("K86").Value = 2332.3125
$ws.Range("L86").Value = 168333.5
$ws.Range("M86").Value = -1209.3125
$ws.Range("N86").Value = -170579.5

$ws.Range("H89").Value = 47605.363
$ws.Range("I89").Value = 2332.3125
$ws.Range("J89").Value = 168333.5
$ws.Range("K89").Value = 11661.5625
$ws.Range("L89").Value = 841667.5
$ws.Range("M89").Value = -6045.5625
$ws.Range("N89").Value = -852899.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 5955056
$ws.Range("I21").Value = 799.5
$ws.Range("J21").Value = 7939808
$ws.Range("K21").Value = 2398.5
$ws.Range("L21").Value = 23819424
$ws.Range("M21").Value = -2225.5
$ws.Range("N21").Value = -23819770

$ws.Range("H34").Value = 498.75
$ws.Range("I34").Value = 231.42857
$ws.Range("J34").Value = 706.6667
$ws.Range("K34").Value = 694.28571
$ws.Range("L34").Value = 2120.0001
$ws.Range("M34").Value = -610.28571
$ws.Range("N34").Value = -2288.0001

$ws.Range("H54").Value = 5800
$ws.Range("J54").Value = 5800
$ws.Range("L54").Value = 17400
$ws.Range("N54").Value = -18518

$ws.Range("H68").Value = 542.5
$ws.Range("I68").Value = 573.3333
$ws.Range("J68").Value = 450
$ws.Range("K68").Value = 1719.9999
$ws.Range("L68").Value = 1350
$ws.Range("M68").Value = -908.9999
$ws.Range("N68").Value = -2972

$ws.Range("H71").Value = 542.5
$ws.Range("I71").Value = 573.3333
$ws.Range("J71").Value = 450
$ws.Range("K71").Value = 5159.9997
$ws.Range("L71").Value = 4050
$ws.Range("M71").Value = -1103.9997
$ws.Range("N71").Value = -12162

$ws.Range("H116").Value = 2666.6667
$ws.Range("I116").Value = 1000
$ws.Range("J116").Value = 6000
$ws.Range("K116").Value = 3000
$ws.Range("L116").Value = 18000
$ws.Range("M116").Value = 442
$ws.Range("N116").Value = -24884

$ws.Range("H117").Value = 63170.125
$ws.Range("J117").Value = 125733
$ws.Range("L117").Value = 377199
$ws.Range("N117").Value = -384083

$ws.Range("H118").Value = 1315
$ws.Range("I118").Value = 1315
$ws.Range("K118").Value = 3945
$ws.Range("M118").Value = -2702

$ws.Range("H119").Value = 3931.0667
$ws.Range("I119").Value = 1269.6364
$ws.Range("J119").Value = 11250
$ws.Range("K119").Value = 3808.9092
$ws.Range("L119").Value = 33750
$ws.Range("M119").Value = 1029.0908
$ws.Range("N119").Value = -43426

$ws.Range("H120").Value = 13316.2
$ws.Range("I120").Value = 9507.5
$ws.Range("J120").Value = 15855.333
$ws.Range("K120").Value = 28522.5
$ws.Range("L120").Value = 47565.999
$ws.Range("M120").Value = -23684.5
$ws.Range("N120").Value = -57241.999

$ws.Range("H121").Value = 26112
$ws.Range("I121").Value = 973.1
$ws.Range("J121").Value = 47061.082
$ws.Range("K121").Value = 2919.3
$ws.Range("L121").Value = 141183.246
$ws.Range("M121").Value = -1609.3
$ws.Range("N121").Value = -143803.246

$ws.Range("H124").Value = 17500
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 17500
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 52500
$ws.Range("M124").ClearContents()
$ws.Range("N124").Value = -62320

$ws.Range("H129").Value = 3125909.8
$ws.Range("J129").Value = 3572382.5
$ws.Range("L129").Value = 10717147.5
$ws.Range("N129").Value = -10727147.5

$ws.Range("H131").Value = 2545.2603
$ws.Range("I131").Value = 460
$ws.Range("J131").Value = 3084.5518
$ws.Range("K131").Value = 1380
$ws.Range("L131").Value = 9253.6554
$ws.Range("M131").Value = 3660
$ws.Range("N131").Value = -19333.6554

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 7121.143
$ws.Range("I13").Value = 341.33334
$ws.Range("J13").Value = 47800
$ws.Range("K13").Value = 341.33334
$ws.Range("L13").Value = 47800
$ws.Range("M13").Value = -202.33334
$ws.Range("N13").Value = -48078

$ws.Range("H80").Value = 3252.8
$ws.Range("I80").Value = 2905.5557
$ws.Range("J80").Value = 3536.9092
$ws.Range("K80").Value = 2905.5557
$ws.Range("L80").Value = 3536.9092
$ws.Range("M80").Value = -1907.5557
$ws.Range("N80").Value = -5532.9092

$ws.Range("H83").Value = 3252.8
$ws.Range("I83").Value = 2905.5557
$ws.Range("J83").Value = 3536.9092
$ws.Range("K83").Value = 14527.7785
$ws.Range("L83").Value = 17684.546
$ws.Range("M83").Value = -9535.7785
$ws.Range("N83").Value = -27668.546

$ws.Range("H135").Value = 52132.223
$ws.Range("J135").Value = 52132.223
$ws.Range("L135").Value = 52132.223
$ws.Range("N135").Value = -62272.223

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2156.5386
$ws.Range("I82").Value = 1114.6666
$ws.Range("J82").Value = 4500.75
$ws.Range("K82").Value = 1114.6666
$ws.Range("L82").Value = 4500.75
$ws.Range("M82").Value = -753.6666
$ws.Range("N82").Value = -5222.75

$ws.Range("H85").Value = 2156.5386
$ws.Range("I85").Value = 1114.6666
$ws.Range("J85").Value = 4500.75
$ws.Range("K85").Value = 1114.6666
$ws.Range("L85").Value = 4500.75
$ws.Range("M85").Value = 133.3334
$ws.Range("N85").Value = -6996.75
